$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 361127.72
$ws.Range("J17").Value = 361127.72
$ws.Range("L17").Value = 1083383.16
$ws.Range("N17").Value = -1083719.16
$ws.Range("H18").Value = 768.8
$ws.Range("I18").Value = 768.8
$ws.Range("K18").Value = 768.8
$ws.Range("M18").Value = -484.8
$ws.Range("H19").Value = 239.76
$ws.Range("I19").Value = 179.66667
$ws.Range("K19").Value = 179.66667
$ws.Range("M19").Value = -4.666670000000011
$ws.Range("H93").Value = 89999
$ws.Range("J93").Value = 89999
$ws.Range("L93").Value = 89999
$ws.Range("N93").Value = -94991
$ws.Range("H98").Value = 1811.6842
$ws.Range("I98").Value = 1538.5
$ws.Range("K98").Value = 1538.5
$ws.Range("M98").Value = -40.5
$ws.Range("H122").Value = 1811.6842
$ws.Range("I122").Value = 1538.5
$ws.Range("K122").Value = 4615.5
$ws.Range("M122").Value = -2165.5
$ws.Range("H132").Value = 31254326
$ws.Range("I132").Value = 35719036
$ws.Range("J132").Value = 1375
$ws.Range("K132").Value = 107157108
$ws.Range("L132").Value = 4125
$ws.Range("M132").Value = -107154578
$ws.Range("N132").Value = -9185
$ws.Range("H135").Value = 3689.1428
$ws.Range("I135").Value = 1590.4
$ws.Range("J135").Value = 8936
$ws.Range("K135").Value = 14313.6
$ws.Range("L135").Value = 80424
$ws.Range("M135").Value = -11778.6
$ws.Range("N135").Value = -85494

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3275.16
$ws.Range("I32").Value = 2819.897
$ws.Range("J32").Value = 17995.334
$ws.Range("K32").Value = 2819.897
$ws.Range("L32").Value = 17995.334
$ws.Range("M32").Value = -2532.897
$ws.Range("N32").Value = -18569.334
$ws.Range("H61").Value = 6390.2036
$ws.Range("I61").Value = 6821.5654
$ws.Range("J61").Value = 3909.875
$ws.Range("K61").Value = 6821.5654
$ws.Range("L61").Value = 3909.875
$ws.Range("M61").Value = -6609.5654
$ws.Range("N61").Value = -4333.875
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
$ws.Range("H103").Value = 90000
$ws.Range("J103").Value = 90000
$ws.Range("L103").Value = 90000
$ws.Range("N103").Value = -92344
$ws.Range("H122").Value = 3291.7354
$ws.Range("I122").Value = 3026.4167
$ws.Range("K122").Value = 9079.250100000001
$ws.Range("M122").Value = -6629.250100000001
$ws.Range("H136").Value = 6390.2036
$ws.Range("I136").Value = 6821.5654
$ws.Range("J136").Value = 3909.875
$ws.Range("K136").Value = 20464.6962
$ws.Range("L136").Value = 11729.625
$ws.Range("M136").Value = -17914.6962
$ws.Range("N136").Value = -16829.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1258.9333
$ws.Range("I64").Value = 486.5
$ws.Range("J64").Value = 1377.7693
$ws.Range("K64").Value = 486.5
$ws.Range("L64").Value = 1377.7693
$ws.Range("M64").Value = -261.5
$ws.Range("N64").Value = -1827.7693
$ws.Range("H67").Value = 1258.9333
$ws.Range("I67").Value = 486.5
$ws.Range("J67").Value = 1377.7693
$ws.Range("K67").Value = 486.5
$ws.Range("L67").Value = 1377.7693
$ws.Range("M67").Value = 293.5
$ws.Range("N67").Value = -2937.7693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 507.69565
$ws.Range("I22").Value = 301.5
$ws.Range("J22").Value = 1250
$ws.Range("K22").Value = 301.5
$ws.Range("L22").Value = 1250
$ws.Range("M22").Value = 48.5
$ws.Range("N22").Value = -1950
$ws.Range("H99").Value = 5614.7617
$ws.Range("I99").Value = 5405
$ws.Range("K99").Value = 5405
$ws.Range("M99").Value = -3907
$ws.Range("H122").Value = 4187.316
$ws.Range("I122").Value = 4303.8887
$ws.Range("K122").Value = 12911.6661
$ws.Range("M122").Value = -10461.6661
$ws.Range("H126").Value = 5614.7617
$ws.Range("I126").Value = 5405
$ws.Range("K126").Value = 16215
$ws.Range("M126").Value = -13745

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 5611.615
$ws.Range("J107").Value = 5611.615
$ws.Range("L107").Value = 16834.845
$ws.Range("N107").Value = -20674.845
$ws.Range("H139").Value = 5558665
$ws.Range("I139").Value = 2078.4167
$ws.Range("J139").Value = 16671839
$ws.Range("K139").Value = 6235.250100000001
$ws.Range("L139").Value = 50015517
$ws.Range("M139").Value = -1095.250100000001
$ws.Range("N139").Value = -50025797
$ws.Range("H141").Value = 23250
$ws.Range("I141").Value = 19666.666
$ws.Range("K141").Value = 58999.99800000001
$ws.Range("M141").Value = -53819.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 45000
$ws.Range("J47").Value = 45000
$ws.Range("L47").Value = 45000
$ws.Range("N47").Value = -46136
$ws.Range("H52").Value = 43297
$ws.Range("J52").Value = 43297
$ws.Range("L52").Value = 43297
$ws.Range("N52").Value = -43815
$ws.Range("H122").Value = 1734.9
$ws.Range("I122").Value = 988.125
$ws.Range("K122").Value = 2964.375
$ws.Range("M122").Value = -514.375
$ws.Range("H132").Value = 3950
$ws.Range("I132").Value = 4070.724
$ws.Range("J132").Value = 2199.5
$ws.Range("K132").Value = 12212.172
$ws.Range("L132").Value = 6598.5
$ws.Range("M132").Value = -9682.172
$ws.Range("N132").Value = -11658.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6617.7085
$ws.Range("I7").Value = 5728.0527
$ws.Range("K7").Value = 5728.0527
$ws.Range("M7").Value = -5616.0527
$ws.Range("H16").Value = 1508.2
$ws.Range("I16").Value = 847.5714
$ws.Range("K16").Value = 847.5714
$ws.Range("M16").Value = -677.5714
$ws.Range("H40").Value = 6777.154
$ws.Range("I40").Value = 6152.8335
$ws.Range("K40").Value = 6152.8335
$ws.Range("M40").Value = -6016.8335
$ws.Range("H82").Value = 11160
$ws.Range("I82").Value = 17933.5
$ws.Range("J82").Value = 999.75
$ws.Range("K82").Value = 17933.5
$ws.Range("L82").Value = 999.75
$ws.Range("M82").Value = -17572.5
$ws.Range("N82").Value = -1721.75
$ws.Range("H85").Value = 11160
$ws.Range("I85").Value = 17933.5
$ws.Range("J85").Value = 999.75
$ws.Range("K85").Value = 17933.5
$ws.Range("L85").Value = 999.75
$ws.Range("M85").Value = -16685.5
$ws.Range("N85").Value = -3495.75
$ws.Range("H126").Value = 6617.7085
$ws.Range("I126").Value = 5728.0527
$ws.Range("K126").Value = 17184.1581
$ws.Range("M126").Value = -14714.1581
$ws.Range("H132").Value = 3418.8408
$ws.Range("I132").Value = 3403.3333
$ws.Range("J132").Value = 3443.4707
$ws.Range("K132").Value = 10209.9999
$ws.Range("L132").Value = 10330.4121
$ws.Range("M132").Value = -7679.999899999999
$ws.Range("N132").Value = -15390.4121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 22505.25
$ws.Range("J39").Value = 29673.666
$ws.Range("L39").Value = 29673.666
$ws.Range("N39").Value = -30499.666
$ws.Range("H43").Value = 25479.5
$ws.Range("J43").Value = 25479.5
$ws.Range("L43").Value = 25479.5
$ws.Range("N43").Value = -25777.5
$ws.Range("H122").Value = 2201.1875
$ws.Range("I122").Value = 1632.8
$ws.Range("J122").Value = 3148.5
$ws.Range("K122").Value = 4898.4
$ws.Range("L122").Value = 9445.5
$ws.Range("M122").Value = -2448.4
$ws.Range("N122").Value = -14345.5
$ws.Range("H126").Value = 3127.3845
$ws.Range("I126").Value = 2050.75
$ws.Range("J126").Value = 3605.889
$ws.Range("K126").Value = 6152.25
$ws.Range("L126").Value = 10817.667
$ws.Range("M126").Value = -3682.25
$ws.Range("N126").Value = -15757.667
$ws.Range("H132").Value = 1313.1945
$ws.Range("I132").Value = 1231.8966
$ws.Range("K132").Value = 3695.6898
$ws.Range("M132").Value = -1165.6898
$ws.Range("H138").Value = 107999.5
$ws.Range("J138").Value = 107999.5
$ws.Range("L138").Value = 107999.5
$ws.Range("N138").Value = -118279.5
